# Apply "Add data for 2022-06-01" update:
#  - Rename the sheet / title text from "...May 23" to "...May 24"
#  - Bump the cumulative carjacking counts for the "through May NN" columns
#    (one column per year: May 2022=B, May 2021=G, May 2020=L, May 2018=V,
#    May 2017=AA) for the neighborhoods that saw a new incident recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet (tab) name and the "May 2022 (through May NN)" header text both
# move from day 23 to day 24.
$ws.Name = "Through 2022-05-24"
$ws.Range("B1").Value = "May 2022 (through May 24)"

# Cell -> new value (only the cells whose counts changed).
$updates = @{
    "G2"  = 7
    "G3"  = 8
    "B4"  = 4
    "L5"  = 5
    "V5"  = 5
    "V6"  = 3
    "G7"  = 2
    "B8"  = 5
    "G8"  = 5
    "V13" = 2
    "V22" = 2
    "G23" = 5
    "AA23" = 3
    "B25" = 5
    "V27" = 1
    "V28" = 2
    "B38" = 3
    "B45" = 4
    "G58" = 2
    "B85" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
